$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(471).Insert()

$ws.Range("A471").Value = 5
$ws.Range("B471").Value = "Macroferia Regional de Talca"
$ws.Range("C471").Value = "Maule"
$ws.Range("D471").Value = 44694
$ws.Range("E471").Value = 7
$ws.Range("F471").Value = 100112004
$ws.Range("G471").Value = "Cebolla"
$ws.Range("H471").Value = "Sin especificar"
$ws.Range("I471").Value = "1a (guarda)"
$ws.Range("J471").Value = 2500
$ws.Range("K471").Value = 7500
$ws.Range("L471").Value = 7500
$ws.Range("M471").Value = 7500
$ws.Range("N471").Value = "`$/malla 25 kilos"
$ws.Range("O471").Value = "Región del Maule"
$ws.Range("P471").Value = 300
$ws.Range("Q471").Value = 25
$ws.Range("R471").Value = "Hortaliza"

Write-Host "done"
